# Generate Report for Handoff
# Adds a new localization-status row (615df727-1b2c-407d-aa24-f1d72a758dc0)
# ahead of the existing ca920b69-... row on all three sheets
# (Overview, zh-cn, de-de), pushing the existing row down to row 3.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0993b519a23d78597bf4f39fadb4c05f7ecf382/e2e/"
$newFile = "615df727-1b2c-407d-aa24-f1d72a758dc0.md"
$oldFile = "ca920b69-3bc7-44c7-85a7-a6ea67819677.md"

# ---------------------------------------------------------------
# Sheet "Overview" (table "Overview", columns A:G)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Rows.Item(2).Insert()

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:G3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "e2e\" + $newFile
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-09-04 18:44:33"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $repoBase + $newFile, "", "", "e2e\" + $newFile)
$ws.Hyperlinks.Add($ws.Range("B3"), $repoBase + $oldFile, "", "", "e2e\" + $oldFile)

# ---------------------------------------------------------------
# Sheet "zh-cn" (table "zh_cn", columns A:P)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Rows.Item(2).Insert()

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "615df727-1b2c-407d-aa24-f1d72a758dc0.e63ae731c01411a3da59f4e4e6da4ac37d918b17.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-04 18:44:29"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = ""
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $repoBase + $newFile, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $repoBase + $oldFile, "", "", $oldFile)

# ---------------------------------------------------------------
# Sheet "de-de" (table "de_de", columns A:P)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Rows.Item(2).Insert()

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P3"))

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "'False"
$ws.Range("G2").Value = "615df727-1b2c-407d-aa24-f1d72a758dc0.e63ae731c01411a3da59f4e4e6da4ac37d918b17.de-de.xlf"
$ws.Range("H2").Value = "2016-09-04 18:44:33"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "'True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "'False"
$ws.Range("P2").Value = ""
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $repoBase + $newFile, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $repoBase + $oldFile, "", "", $oldFile)
